$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the species/observation data between row 26 and row 27
# (columns A, B, E, F, G, H, Q, R, Z, AB) while leaving the other
# shared columns (which are identical in both rows) untouched.

$cols = @("A","B","E","F","G","H","Q","R","Z","AB")

foreach ($col in $cols) {
    $r26 = $ws.Range($col + "26")
    $r27 = $ws.Range($col + "27")
    $tmp = $r26.Value2
    $r26.Value = $r27.Value2
    $r27.Value = $tmp
}
